{"js": "// Remove the \"Floor texture\" and \"Wallpaper texture\" asset-reference list\n// items (each a bullet with a trailing hyperlink) from under the\n// \"Asset References\" heading. This matches the diff: the two bulleted\n// <w:p> paragraphs are deleted in their entirety, leaving the heading\n// immediately followed by the (already-empty) final list item.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"text\");\nawait context.sync();\n\n// Collect the paragraphs whose text begins with the asset labels that are\n// being removed; gather first, then delete, so the live collection isn't\n// mutated mid-iteration.\nconst paragraphsToRemove = [];\nfor (const paragraph of body.paragraphs.items) {\n  const text = paragraph.text.trim();\n  if (text.indexOf(\"Floor texture\") === 0 || text.indexOf(\"Wallpaper texture\") === 0) {\n    paragraphsToRemove.push(paragraph);\n  }\n}\n\nfor (const paragraph of paragraphsToRemove) {\n  paragraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Floor texture\" and \"Wallpaper texture\" asset-reference list\n# items (each a bullet with a trailing hyperlink) from under the\n# \"Asset References\" heading. This matches the diff: the two bulleted\n# paragraphs are deleted in their entirety, leaving the heading\n# immediately followed by the (already-empty) final list item.\n\n$d = $word.ActiveDocument\n$wdParagraph = 4\n\n$labels = @(\"Floor texture\", \"Wallpaper texture\")\n\nforeach ($label in $labels) {\n    $range = $d.Content\n    $found = $range.Find.Execute($label)\n    if ($found) {\n        [void]$range.Expand($wdParagraph)\n        $range.Delete()\n    }\n}\n"}
